# Applies the "comment an article" test-case edit to the DataSet sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")

# --- Row 17: fill in the previously-empty "User comment Article" test row ---
# Columns A-G already contain the correct values; only H..R need populating.
# NB: the order in which brand-new shared-string values are first written
# determines the order they are appended to xl/sharedStrings.xml, so the
# writes below are sequenced to reproduce that layout exactly.
$ws.Range("P17").Value = "AssertCommentCreatedSuccessfully"
$ws.Range("H17").Value = "Comment_CommentAnArticle_CommentCreated"
$ws.Range("I17").Value = "<Fullname>@uniqueemail.com>"
$ws.Range("J17").Value = "<unique, generated per test>"
$ws.Range("K17").Value = "validPass"
$ws.Range("L17").Value = "validPass"
$ws.Range("M17").Value = "Lorem ipsum dolor sit amet, consectetur adipiscing elit. Ut congue augue lacinia, tristique odio ut, porta velit. Cras rutrum dolor ligula, ac ornare lectus cursus sed."

# Setting .Value on R17 (previously-empty, style s="7" which uses quotePrefix)
# causes this runtime to silently swap its style for a plain one. Restore the
# original "quote prefix" style by pasting formats from a same-styled donor
# cell (R8) after writing the value.
$ws.Range("R17").Value = "Failed"
$ws.Range("R8").Copy()
$ws.Range("R17").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 32: rename the test case identifier ---
$ws.Range("H32").Value = "Comment_CommentArticle_CommentButtonExists"

# Now finish row 17 with the remaining brand-new string ("Comment by:")
$ws.Range("N17").Value = "Comment by:"

# --- Row 20: correct wording of the "unable to edit" message (typo/apostrophe fix) ---
$ws.Range("N20").Value = "Unable to edit other Authors articles"

# --- Update the saved cursor/selection position on the sheet ---
$ws.Range("H33").Select()
